# Slide 25 - "Content Placeholder 2": remove the blank paragraph and the
# "    setType(Type.Integer);" paragraph that used to sit between the
# "      }" (end of catch block) paragraph and the final "  }" paragraph.
$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(25)
$shp = $s.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

# Locate the text to remove: the paragraph mark that ends the "      }"
# paragraph, the (empty) blank-line paragraph, the
# "    setType(Type.Integer);" paragraph and its own paragraph mark.
$cr     = [char]13
$needle = "$cr$cr    setType(Type.Integer);$cr"
$full   = $tr.Text
$start  = $full.IndexOf($needle) + 1   # COM ranges are 1-based

$tr.Characters($start, $needle.Length).Delete()
